$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (45189 -> 45190, i.e. 2023-09-20 -> 2023-09-21) for every data row (rows 2-372).
$ws.Range("C2:C372").Value = 45190
